$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.954.29"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.235.08"
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.25"
$ws.Range("E5").Value = "  -3.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.30"
$ws.Range("E6").Value = "  +6.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0799"
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.62"
$ws.Range("E11").Value = "  +7.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.09"
$ws.Range("E12").Value = "  -9.64%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.40"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "2.577.40"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.17"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "2.239.27"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.728"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "39.861.91"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").Value = "0.0₃0894"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.80"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.62"
$ws.Range("E22").Value = "  +5.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.54"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.01"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  +4.30%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.97"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.27"
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.27"
$ws.Range("E31").Value = "  +7.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.11"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0710"
$ws.Range("E35").Value = "  +2.07%  "
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.52"
$ws.Range("E37").Value = "  +10.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").Value = "  +4.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.79"
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("D43").Value = "1.962.94"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0271"
$ws.Range("E45").Value = "  +6.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.59"
$ws.Range("E46").Value = "  +5.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.21"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").Value = "2.448.79"
$ws.Range("E49").Value = "  -3.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.97"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("E51").Value = "  +11.68%  "
